# Fixed error in tp/sl for HA_VWAP that created astronomical gains.
#
# - Row 2 (Test #1, BTCUSDT) and Row 3 (Test #2, BTCUSDT) for the HA_VWAP
#   strategy had their settings corrected from DistVWAP_PCT 0.0 to 0.05,
#   and now use NB_SIGNALS 1 and 2 respectively.
# - New test rows were added (4-7) extending HA_VWAP coverage across
#   BTCUSDT/ETHUSDT with NB_SIGNALS 1, 2 and 3, all using DistVWAP_PCT 0.05.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Correct the existing two HA_VWAP test rows ---
$ws.Range("K2").Value = '{"EMA": 200, "DistVWAP_PCT": 0.05, "NB_SIGNALS": 1}'
$ws.Range("K3").Value = '{"EMA": 200, "DistVWAP_PCT": 0.05, "NB_SIGNALS": 2}'

# --- Row 4: Test #3, BTCUSDT, NB_SIGNALS 3 ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Bybit"
$ws.Range("C4").Value = "BTCUSDT"
$ws.Range("D4").Value = 44562
$ws.Range("E4").Value = 44926
$ws.Range("F4").Value = "5m"
$ws.Range("G4").Value = 7
$ws.Range("H4").Value = 7
$ws.Range("I4").Value = "HA_VWAP"
$ws.Range("J4").Value = "VWAP_Touch"
$ws.Range("K4").Value = '{"EMA": 200, "DistVWAP_PCT": 0.05, "NB_SIGNALS": 3}'

# --- Row 5: Test #4, ETHUSDT, NB_SIGNALS 1 ---
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Bybit"
$ws.Range("C5").Value = "ETHUSDT"
$ws.Range("D5").Value = 44562
$ws.Range("E5").Value = 44926
$ws.Range("F5").Value = "5m"
$ws.Range("G5").Value = 7
$ws.Range("H5").Value = 7
$ws.Range("I5").Value = "HA_VWAP"
$ws.Range("J5").Value = "VWAP_Touch"
$ws.Range("K5").Value = '{"EMA": 200, "DistVWAP_PCT": 0.05, "NB_SIGNALS": 1}'

# --- Row 6: Test #5, ETHUSDT, NB_SIGNALS 2 ---
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Bybit"
$ws.Range("C6").Value = "ETHUSDT"
$ws.Range("D6").Value = 44562
$ws.Range("E6").Value = 44926
$ws.Range("F6").Value = "5m"
$ws.Range("G6").Value = 7
$ws.Range("H6").Value = 7
$ws.Range("I6").Value = "HA_VWAP"
$ws.Range("J6").Value = "VWAP_Touch"
$ws.Range("K6").Value = '{"EMA": 200, "DistVWAP_PCT": 0.05, "NB_SIGNALS": 2}'

# --- Row 7: Test #6, ETHUSDT, NB_SIGNALS 3 ---
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Bybit"
$ws.Range("C7").Value = "ETHUSDT"
$ws.Range("D7").Value = 44562
$ws.Range("E7").Value = 44926
$ws.Range("F7").Value = "5m"
$ws.Range("G7").Value = 7
$ws.Range("H7").Value = 7
$ws.Range("I7").Value = "HA_VWAP"
$ws.Range("J7").Value = "VWAP_Touch"
$ws.Range("K7").Value = '{"EMA": 200, "DistVWAP_PCT": 0.05, "NB_SIGNALS": 3}'

# Rows 6 and 7 reused cells (B/C/F/I) that already carried stray per-cell
# formatting from the old placeholder rows, so their number formats need to
# be re-stamped to match the rest of the table (rows 4 & 5 are already
# correct and serve as the format source).
$ws.Range("A5:K5").Copy()
$ws.Range("A6:K6").PasteSpecial(-4122)
$ws.Range("A5:K5").Copy()
$ws.Range("A7:K7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Restore the view: Sheet1 active, cell C10 selected ---
$ws.Activate()
$ws.Range("C10").Select()
